{"js": "// Update copyright years in the D4AS copyright notice:\n//   \"Advance Steel 2023\" -> \"Advance Steel 2024\"\n//   \"\u00a9 2022 Autodesk, Inc.\" -> \"\u00a9 2023 Autodesk, Inc.\"\nconst body = context.document.body;\n\n// --- \"Steel 2023\" -> \"Steel 2024\" --------------------------------------\nconst steelResults = body.search(\"Steel 2023\", { matchCase: true, matchWholeWord: false });\nsteelResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < steelResults.items.length; i++) {\n  // Narrow to just the trailing \"3\" so the edit mirrors retyping the\n  // last digit of the year rather than rewriting the whole run.\n  const yearDigit = steelResults.items[i].search(\"3\", { matchCase: true });\n  yearDigit.load(\"items\");\n  await context.sync();\n\n  if (yearDigit.items.length > 0) {\n    yearDigit.items[yearDigit.items.length - 1].insertText(\"4\", Word.InsertLocation.replace);\n  } else {\n    steelResults.items[i].insertText(\"Steel 2024\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// --- \"2022\" -> \"2023\" (the Autodesk, Inc. copyright year) ---------------\nconst yearResults = body.search(\"2022\", { matchCase: true, matchWholeWord: true });\nyearResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < yearResults.items.length; i++) {\n  const lastDigit = yearResults.items[i].search(\"2\", { matchCase: true });\n  lastDigit.load(\"items\");\n  await context.sync();\n\n  if (lastDigit.items.length > 0) {\n    lastDigit.items[lastDigit.items.length - 1].insertText(\"3\", Word.InsertLocation.replace);\n  } else {\n    yearResults.items[i].insertText(\"2023\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update copyright years in the D4AS copyright notice:\n#   \"Advance Steel 2023\" -> \"Advance Steel 2024\"\n#   \"(c) 2022 Autodesk, Inc.\" -> \"(c) 2023 Autodesk, Inc.\"\n$d = $word.ActiveDocument\n\n# --- \"Steel 2023\" -> \"Steel 2024\" ---------------------------------------\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.Text = \"Steel 2023\"\n$find1.Find.MatchCase = $true\n$find1.Find.MatchWholeWord = $false\n$found1 = $find1.Find.Execute()\n\nif ($found1) {\n  # Narrow to just the trailing \"3\" so the edit mirrors retyping the last\n  # digit of the year rather than rewriting the whole phrase.\n  $yearDigit = $d.Range($find1.End - 1, $find1.End)\n  $yearDigit.Text = \"4\"\n}\n\n# --- \"2022\" -> \"2023\" (the Autodesk, Inc. copyright year) ---------------\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.Text = \"2022\"\n$find2.Find.MatchCase = $true\n$find2.Find.MatchWholeWord = $true\n$found2 = $find2.Find.Execute()\n\nif ($found2) {\n  $lastDigit = $d.Range($find2.End - 1, $find2.End)\n  $lastDigit.Text = \"3\"\n}\n"}
